{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1) Insert a new subscript \"=\" run at the very start of the document,\n//    right before \"Effective Figures\".\nconst eqRange = body.insertText(\"=\", Word.InsertLocation.start);\neqRange.font.subscript = true;\n\n// 2) & 3) Move the \"_GoBack\" bookmark from around the two figures to\n//    just before \"Are the numbers and labels on each axis...\" text,\n//    which splits that run into \"  \" + \"Are the numbers...\".\ndoc.deleteBookmark(\"_GoBack\");\n\nconst results = body.search(\"Are the numbers and labels\", { matchCase: true });\nawait context.sync();\n\nconst target = results.items[0];\nconst startRange = target.getRange(Word.RangeLocation.start);\nstartRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new subscript \"=\" run at the very start of the document,\n#    right before \"Effective Figures\".\n$startRng = $d.Range(0, 0)\n$startRng.InsertBefore(\"=\")\n$eqRng = $d.Range(0, 1)\n$eqRng.Font.Subscript = $true\n\n# 2) & 3) Move the \"_GoBack\" bookmark from around the two figures to\n#    just before \"Are the numbers and labels on each axis...\" text,\n#    splitting that run into \"  \" + \"Are the numbers...\".\n$findRng = $d.Content\n$findRng.Find.Execute(\"Are the numbers and labels\")\n$bmRange = $d.Range($findRng.Start, $findRng.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
